# Last Order column added to AmazonOrder
# Rebuilds the data rows of Sheet1 to reflect the new, longer order list.
# Row 1 (header: sku / ORD / DESCRIPTION) is unchanged.
# Rows 2-10 previously held 9 data rows; the sheet now holds 29 data rows
# (rows 2-30), so the existing style used on row 2 (style index "2" in the
# source file) is copied down across the whole new range before values are
# written, ensuring the newly added rows pick up the same formatting as the
# original data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2 through 30: sku, ORD, DESCRIPTION
$data = @(
    @(2, 'AL_EASY BRAID 18_34', 2, '14657 34'),
    @(3, 'AL_NAOMI_T14/88', 1, '10397 T14/88'),
    @(4, 'AL_SASHA_1B', 1, '10521 1B'),
    @(5, 'BY_AUSTIN_SP427', 1, 'AUSTIN'),
    @(6, 'BY_DOLLY_2', 1, 'DOLLY'),
    @(7, 'BY_KOBY_T27B', 2, 'KOBY'),
    @(8, 'BY_MALIA_2', 1, 'MALIA'),
    @(9, 'BY_MATA_SP427', 1, 'MATA'),
    @(10, 'BY_PEDY_SP427', 1, 'PEDY'),
    @(11, 'BY_STB COGIC_M27/613', 1, 'STB COGIC'),
    @(12, 'BY_T2HB DENON_SP1B/30', 1, 'T2HB DENON'),
    @(13, 'HZ_EXOTIC_1B', 1, 'IDPEXO'),
    @(14, 'HZ_HAYLIE_1B', 1, 'IBBHAY'),
    @(15, 'HZ_IDPWKS24_2', 2, 'IDPWKS24IDPWKS24'),
    @(16, 'HZ_MILEY_1B', 1, 'HEFMIL'),
    @(17, 'HZ_SB2XB48_3T2/30/27', 2, 'SB2XB48'),
    @(18, 'HZ_VOGUE CROP_1B', 1, 'IHBVCR'),
    @(19, 'OT_ PQWPNBF42_DR2T1B/2730', 1, 'PQWPNBF42'),
    @(20, 'OT_3C WHIRLY LOOP_2', 10, 'KXBIW'),
    @(21, 'OT_DUVESSA REMI YAKI 16_LT1B/433', 4, 'HWDVY16'),
    @(22, 'OT_GOLD OCEAN BODY 16-18-20_NBLK', 1, 'HWMGO161820'),
    @(23, 'OT_LUXELINE NATURAL BODY 10_NBLK', 2, 'HWLLNB10'),
    @(24, 'OT_ROLL UP 44PCS_280', 1, 'HWVBRU234'),
    @(25, 'OT_ROLL UP 44PCS_33', 1, 'HWVBRU234'),
    @(26, 'OT_ROLL UP 44PCS_4', 1, 'HWVBRU234'),
    @(27, 'OT_SPIRALLY_4', 1, 'HWPBSP5'),
    @(28, 'OT_TOYA_1', 2, 'QPNTOYQPNTOY'),
    @(29, 'VF_BRIE-V_99J', 1, '45614'),
    @(30, 'VF_ELSIE_613', 1, '46700')
)

# Copy the formatting of the existing data row (row 2) down to every row
# that will now hold data, including the newly appended rows 11-30.
$srcRow = $ws.Range("A2:C2")
$destRows = $ws.Range("A3:C30")
$srcRow.Copy($destRows)

# DESCRIPTION values that are purely numeric digits (e.g. "45614", "46700",
# on the new rows 29-30) must stay text, matching the original inlineStr
# cell type, instead of being auto-converted to numbers. Force a text
# number format just on those cells before writing their values.
$ws.Range("C29:C30").NumberFormat = "@"

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
}
